$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ntn1"
$ws.Cells.Item(2,3).Value = "Mcam"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.990837
$ws.Cells.Item(2,8).Value = 5.972511000000001
$ws.Cells.Item(2,9).Value = 0.1122845585713437
$ws.Cells.Item(2,10).Value = 0.1122845585713437
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 87.038094
$ws.Cells.Item(2,14).Value = 261.114282
$ws.Cells.Item(2,15).Value = 0.7371871251537216
$ws.Cells.Item(2,16).Value = 0.7371871251537216
$ws.Cells.Item(2,17).Value = 173.278657944678
$ws.Cells.Item(2,18).Value = 1559.507921502102
$ws.Cells.Item(2,19).Value = 0.08277473093236357
$ws.Cells.Item(2,20).Value = 0.08277473093236357

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ntn1"
$ws.Cells.Item(3,3).Value = "Mcam"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.990837
$ws.Cells.Item(3,8).Value = 5.972511000000001
$ws.Cells.Item(3,9).Value = 0.1122845585713437
$ws.Cells.Item(3,10).Value = 0.1122845585713437
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.255465
$ws.Cells.Item(3,14).Value = 6.766394999999999
$ws.Cells.Item(3,15).Value = 0.01910312695076754
$ws.Cells.Item(3,16).Value = 0.01910312695076754
$ws.Cells.Item(3,17).Value = 4.490263174204999
$ws.Cells.Item(3,18).Value = 40.412368567845
$ws.Cells.Item(3,19).Value = 0.002144986176999273
$ws.Cells.Item(3,20).Value = 0.002144986176999274

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ntn1"
$ws.Cells.Item(4,3).Value = "Mcam"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.990837
$ws.Cells.Item(4,8).Value = 5.972511000000001
$ws.Cells.Item(4,9).Value = 0.1122845585713437
$ws.Cells.Item(4,10).Value = 0.1122845585713437
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.5269253333333334
$ws.Cells.Item(4,14).Value = 1.580776
$ws.Cells.Item(4,15).Value = 0.004462903009464643
$ws.Cells.Item(4,16).Value = 0.004462903009464643
$ws.Cells.Item(4,17).Value = 1.049022449837334
$ws.Cells.Item(4,18).Value = 9.441202048536001
$ws.Cells.Item(4,19).Value = 0.000501115094364459
$ws.Cells.Item(4,20).Value = 0.000501115094364459

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Ntn1"
$ws.Cells.Item(5,3).Value = "Mcam"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.990837
$ws.Cells.Item(5,8).Value = 5.972511000000001
$ws.Cells.Item(5,9).Value = 0.1122845585713437
$ws.Cells.Item(5,10).Value = 0.1122845585713437
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 28.247359
$ws.Cells.Item(5,14).Value = 84.74207700000001
$ws.Cells.Item(5,15).Value = 0.2392468448860462
$ws.Cells.Item(5,16).Value = 0.2392468448860462
$ws.Cells.Item(5,17).Value = 56.23588744948301
$ws.Cells.Item(5,18).Value = 506.1229870453471
$ws.Cells.Item(5,19).Value = 0.02686372636761644
$ws.Cells.Item(5,20).Value = 0.02686372636761645

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ntn1"
$ws.Cells.Item(6,3).Value = "Mcam"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 11.42765333333333
$ws.Cells.Item(6,8).Value = 34.28296
$ws.Cells.Item(6,9).Value = 0.6445274073365515
$ws.Cells.Item(6,10).Value = 0.6445274073365515
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 87.038094
$ws.Cells.Item(6,14).Value = 261.114282
$ws.Cells.Item(6,15).Value = 0.7371871251537216
$ws.Cells.Item(6,16).Value = 0.7371871251537216
$ws.Cells.Item(6,17).Value = 994.64116502608
$ws.Cells.Item(6,18).Value = 8951.770485234721
$ws.Cells.Item(6,19).Value = 0.4751373064972141
$ws.Cells.Item(6,20).Value = 0.4751373064972141

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ntn1"
$ws.Cells.Item(7,3).Value = "Mcam"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 11.42765333333333
$ws.Cells.Item(7,8).Value = 34.28296
$ws.Cells.Item(7,9).Value = 0.6445274073365515
$ws.Cells.Item(7,10).Value = 0.6445274073365515
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.255465
$ws.Cells.Item(7,14).Value = 6.766394999999999
$ws.Cells.Item(7,15).Value = 0.01910312695076754
$ws.Cells.Item(7,16).Value = 0.01910312695076754
$ws.Cells.Item(7,17).Value = 25.77467212546666
$ws.Cells.Item(7,18).Value = 231.9720491292
$ws.Cells.Item(7,19).Value = 0.0123124888855992
$ws.Cells.Item(7,20).Value = 0.01231248888559921

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Ntn1"
$ws.Cells.Item(8,3).Value = "Mcam"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 11.42765333333333
$ws.Cells.Item(8,8).Value = 34.28296
$ws.Cells.Item(8,9).Value = 0.6445274073365515
$ws.Cells.Item(8,10).Value = 0.6445274073365515
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.5269253333333334
$ws.Cells.Item(8,14).Value = 1.580776
$ws.Cells.Item(8,15).Value = 0.004462903009464643
$ws.Cells.Item(8,16).Value = 0.004462903009464643
$ws.Cells.Item(8,17).Value = 6.021520041884445
$ws.Cells.Item(8,18).Value = 54.19368037696
$ws.Cells.Item(8,19).Value = 0.00287646330588474
$ws.Cells.Item(8,20).Value = 0.00287646330588474

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Ntn1"
$ws.Cells.Item(9,3).Value = "Mcam"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 11.42765333333333
$ws.Cells.Item(9,8).Value = 34.28296
$ws.Cells.Item(9,9).Value = 0.6445274073365515
$ws.Cells.Item(9,10).Value = 0.6445274073365515
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 28.247359
$ws.Cells.Item(9,14).Value = 84.74207700000001
$ws.Cells.Item(9,15).Value = 0.2392468448860462
$ws.Cells.Item(9,16).Value = 0.2392468448860462
$ws.Cells.Item(9,17).Value = 322.8010262342134
$ws.Cells.Item(9,18).Value = 2905.20923610792
$ws.Cells.Item(9,19).Value = 0.1542011486478534
$ws.Cells.Item(9,20).Value = 0.1542011486478534

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Ntn1"
$ws.Cells.Item(10,3).Value = "Mcam"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.2266433333333333
$ws.Cells.Item(10,8).Value = 0.67993
$ws.Cells.Item(10,9).Value = 0.01278283789002879
$ws.Cells.Item(10,10).Value = 0.01278283789002879
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 87.038094
$ws.Cells.Item(10,14).Value = 261.114282
$ws.Cells.Item(10,15).Value = 0.7371871251537216
$ws.Cells.Item(10,16).Value = 0.7371871251537216
$ws.Cells.Item(10,17).Value = 19.72660375114
$ws.Cells.Item(10,18).Value = 177.53943376026
$ws.Cells.Item(10,19).Value = 0.009423343515456388
$ws.Cells.Item(10,20).Value = 0.009423343515456388

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Ntn1"
$ws.Cells.Item(11,3).Value = "Mcam"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.2266433333333333
$ws.Cells.Item(11,8).Value = 0.67993
$ws.Cells.Item(11,9).Value = 0.01278283789002879
$ws.Cells.Item(11,10).Value = 0.01278283789002879
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 2.255465
$ws.Cells.Item(11,14).Value = 6.766394999999999
$ws.Cells.Item(11,15).Value = 0.01910312695076754
$ws.Cells.Item(11,16).Value = 0.01910312695076754
$ws.Cells.Item(11,17).Value = 0.5111861058166666
$ws.Cells.Item(11,18).Value = 4.600674952349999
$ws.Cells.Item(11,19).Value = 0.0002441921750043014
$ws.Cells.Item(11,20).Value = 0.0002441921750043015

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Ntn1"
$ws.Cells.Item(12,3).Value = "Mcam"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.2266433333333333
$ws.Cells.Item(12,8).Value = 0.67993
$ws.Cells.Item(12,9).Value = 0.01278283789002879
$ws.Cells.Item(12,10).Value = 0.01278283789002879
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.5269253333333334
$ws.Cells.Item(12,14).Value = 1.580776
$ws.Cells.Item(12,15).Value = 0.004462903009464643
$ws.Cells.Item(12,16).Value = 0.004462903009464643
$ws.Cells.Item(12,17).Value = 0.1194241139644445
$ws.Cells.Item(12,18).Value = 1.07481702568
$ws.Cells.Item(12,19).Value = 0.00005704856568890816
$ws.Cells.Item(12,20).Value = 0.00005704856568890816

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Ntn1"
$ws.Cells.Item(13,3).Value = "Mcam"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.2266433333333333
$ws.Cells.Item(13,8).Value = 0.67993
$ws.Cells.Item(13,9).Value = 0.01278283789002879
$ws.Cells.Item(13,10).Value = 0.01278283789002879
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 28.247359
$ws.Cells.Item(13,14).Value = 84.74207700000001
$ws.Cells.Item(13,15).Value = 0.2392468448860462
$ws.Cells.Item(13,16).Value = 0.2392468448860462
$ws.Cells.Item(13,17).Value = 6.402075601623334
$ws.Cells.Item(13,18).Value = 57.61868041461001
$ws.Cells.Item(13,19).Value = 0.003058253633879192
$ws.Cells.Item(13,20).Value = 0.003058253633879192

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Ntn1"
$ws.Cells.Item(14,3).Value = "Mcam"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 4.085149333333333
$ws.Cells.Item(14,8).Value = 12.255448
$ws.Cells.Item(14,9).Value = 0.230405196202076
$ws.Cells.Item(14,10).Value = 0.230405196202076
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 87.038094
$ws.Cells.Item(14,14).Value = 261.114282
$ws.Cells.Item(14,15).Value = 0.7371871251537216
$ws.Cells.Item(14,16).Value = 0.7371871251537216
$ws.Cells.Item(14,17).Value = 355.563611678704
$ws.Cells.Item(14,18).Value = 3200.072505108336
$ws.Cells.Item(14,19).Value = 0.1698517442086876
$ws.Cells.Item(14,20).Value = 0.1698517442086876

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Ntn1"
$ws.Cells.Item(15,3).Value = "Mcam"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 4.085149333333333
$ws.Cells.Item(15,8).Value = 12.255448
$ws.Cells.Item(15,9).Value = 0.230405196202076
$ws.Cells.Item(15,10).Value = 0.230405196202076
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 2.255465
$ws.Cells.Item(15,14).Value = 6.766394999999999
$ws.Cells.Item(15,15).Value = 0.01910312695076754
$ws.Cells.Item(15,16).Value = 0.01910312695076754
$ws.Cells.Item(15,17).Value = 9.213911341106664
$ws.Cells.Item(15,18).Value = 82.92520206995999
$ws.Cells.Item(15,19).Value = 0.004401459713164761
$ws.Cells.Item(15,20).Value = 0.004401459713164762

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Ntn1"
$ws.Cells.Item(16,3).Value = "Mcam"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 4.085149333333333
$ws.Cells.Item(16,8).Value = 12.255448
$ws.Cells.Item(16,9).Value = 0.230405196202076
$ws.Cells.Item(16,10).Value = 0.230405196202076
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.5269253333333334
$ws.Cells.Item(16,14).Value = 1.580776
$ws.Cells.Item(16,15).Value = 0.004462903009464643
$ws.Cells.Item(16,16).Value = 0.004462903009464643
$ws.Cells.Item(16,17).Value = 2.152568674183111
$ws.Cells.Item(16,18).Value = 19.373118067648
$ws.Cells.Item(16,19).Value = 0.001028276043526537
$ws.Cells.Item(16,20).Value = 0.001028276043526537

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Ntn1"
$ws.Cells.Item(17,3).Value = "Mcam"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 4.085149333333333
$ws.Cells.Item(17,8).Value = 12.255448
$ws.Cells.Item(17,9).Value = 0.230405196202076
$ws.Cells.Item(17,10).Value = 0.230405196202076
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 28.247359
$ws.Cells.Item(17,14).Value = 84.74207700000001
$ws.Cells.Item(17,15).Value = 0.2392468448860462
$ws.Cells.Item(17,16).Value = 0.2392468448860462
$ws.Cells.Item(17,17).Value = 115.3946797872773
$ws.Cells.Item(17,18).Value = 1038.552118085496
$ws.Cells.Item(17,19).Value = 0.05512371623669711
$ws.Cells.Item(17,20).Value = 0.05512371623669712
